$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the descriptive text cells to reference FY 2012-2016 instead of FY 2011-2016
$ws.Range("A3").Value = "This table shows the grant awards and award dollars AHRQ made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the AHRQ page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars AHRQ made for FY 2012-2016."
